$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: the "sin fechas..." long observation note is removed from F6.
$ws.Range("F6").Value = ""

# Row 10: add delivery date and resource code REC10.
$ws.Range("E10").Value = 42085
$ws.Range("F10").Value = "REC10"

# Row 11: add delivery date and resource code REC160.
$ws.Range("E11").Value = 42086
$ws.Range("F11").Value = "REC160"

# Row 6 height shrinks now that the long note text is gone.
$ws.Rows(6).RowHeight = 30.75

# Update the active selection to reflect where the editor left off.
$ws.Range("E11").Select()
